# Results from stan models 2016-03-30.xlsx
# "Format, streamline analysis, order work
#  also changed traits back to day, not effect size"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window position (best effort; engine may not persist this) ---
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 16300
    $win.Top = 880
} catch {}

# --- Column A width (20.5 chars in the saved XML) ---
$ws.Columns.Item(1).ColumnWidth = 19.666666666666668

# --- New results table: "chilling as two levels, with interactions" / leafout ---
# Body of the new table pasted first (this is what creates the new shared
# strings mu_b_site..mu_b_inter_wc2 in that order), then the header rows are
# filled in above it.

$ws.Range("A26").Value = "mu_b_site"
$ws.Range("B26").Value = 214
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 0.1
$ws.Range("F26").Value = 2.1
$ws.Range("G26").Value = -1.1
$ws.Range("H26").Value = 1.7
$ws.Range("I26").Value = 3
$ws.Range("J26").Value = 4.4
$ws.Range("K26").Value = 6.8

$ws.Range("A27").Value = "mu_b_inter_wp"
$ws.Range("B27").Value = 320
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 3.5
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0.8
$ws.Range("G27").Value = 2.1
$ws.Range("H27").Value = 3
$ws.Range("I27").Value = 3.5
$ws.Range("J27").Value = 4.1
$ws.Range("K27").Value = 5.1

$ws.Range("A28").Value = "mu_b_inter_ws"
$ws.Range("B28").Value = 141
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = -0.6
$ws.Range("E28").Value = 0.1
$ws.Range("F28").Value = 0.8
$ws.Range("G28").Value = -2.4
$ws.Range("H28").Value = -1.1
$ws.Range("I28").Value = -0.6
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 1

$ws.Range("A29").Value = "mu_b_inter_ps"
$ws.Range("B29").Value = 217
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = -1
$ws.Range("E29").Value = 0.1
$ws.Range("F29").Value = 0.8
$ws.Range("G29").Value = -2.6
$ws.Range("H29").Value = -1.6
$ws.Range("I29").Value = -1
$ws.Range("J29").Value = -0.4
$ws.Range("K29").Value = 0.7

$ws.Range("A30").Value = "mu_b_inter_wc1"
$ws.Range("B30").Value = 769
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 10.2
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 1.2
$ws.Range("G30").Value = 7.8
$ws.Range("H30").Value = 9.5
$ws.Range("I30").Value = 10.1
$ws.Range("J30").Value = 10.9
$ws.Range("K30").Value = 12.5

$ws.Range("A31").Value = "mu_b_inter_wc2"
$ws.Range("B31").Value = 224
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 11.3
$ws.Range("E31").Value = 0.1
$ws.Range("F31").Value = 1.2
$ws.Range("G31").Value = 9
$ws.Range("H31").Value = 10.4
$ws.Range("I31").Value = 11.3
$ws.Range("J31").Value = 12.1
$ws.Range("K31").Value = 13.7

# Header rows for the new table
$ws.Range("A20").Value = "lday_site_sp_chill_inter"
$ws.Range("B20").Value = "chilling as two levels, with interactions"

$ws.Range("B21").Value = "n_eff"
$ws.Range("C21").Value = "Rhat"
$ws.Range("D21").Value = "mean"
$ws.Range("E21").Value = "mcse"
$ws.Range("F21").Value = "sd"
$ws.Range("G21").NumberFormat = "0.00%"
$ws.Range("G21").Value = 0.025
$ws.Range("H21").NumberFormat = "0%"
$ws.Range("H21").Value = 0.25
$ws.Range("I21").NumberFormat = "0%"
$ws.Range("I21").Value = 0.5
$ws.Range("J21").NumberFormat = "0%"
$ws.Range("J21").Value = 0.75
$ws.Range("K21").NumberFormat = "0.00%"
$ws.Range("K21").Value = 0.975

# Re-insert the existing rows for mu_b_warm..mu_b_photo in the new table
$ws.Range("A22").Value = "mu_b_warm"
$ws.Range("B22").Value = 236
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = -21.9
$ws.Range("E22").Value = 0.1
$ws.Range("F22").Value = 1.7
$ws.Range("G22").Value = -25.3
$ws.Range("H22").Value = -23
$ws.Range("I22").Value = -21.9
$ws.Range("J22").Value = -20.7
$ws.Range("K22").Value = -18.6

$ws.Range("A23").Value = "mu_b_chill1"
$ws.Range("B23").Value = 475
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = -26.4
$ws.Range("E23").Value = 0.1
$ws.Range("F23").Value = 3.1
$ws.Range("G23").Value = -32.4
$ws.Range("H23").Value = -28.4
$ws.Range("I23").Value = -26.4
$ws.Range("J23").Value = -24.4
$ws.Range("K23").Value = -20.3

$ws.Range("A24").Value = "mu_b_chill2"
$ws.Range("B24").Value = 345
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = -26.1
$ws.Range("E24").Value = 0.2
$ws.Range("F24").Value = 3.1
$ws.Range("G24").Value = -32.1
$ws.Range("H24").Value = -28.3
$ws.Range("I24").Value = -26.2
$ws.Range("J24").Value = -24
$ws.Range("K24").Value = -20.2

$ws.Range("A25").Value = "mu_b_photo"
$ws.Range("B25").Value = 270
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = -13.7
$ws.Range("E25").Value = 0.1
$ws.Range("F25").Value = 1.7
$ws.Range("G25").Value = -17
$ws.Range("H25").Value = -14.8
$ws.Range("I25").Value = -13.7
$ws.Range("J25").Value = -12.6
$ws.Range("K25").Value = -10.4

# Title row added above the table last
$ws.Range("A19").Value = "leafout"

# Selection ends on F22, matching the final cursor position in the workbook
[void]$ws.Range("F22").Select()
